$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-08"

# Update the header label cell (I1) that shows the "through" date
$ws.Range("I1").Value = "2022 (through 06-08)"

# Update the data values for the newly added day (June row / Total row)
$ws.Range("I7").Value = 29
$ws.Range("I14").Value = 692
